$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "57.970.48"
$ws.Range("E2").Value = "  -0.71%  "
$ws.Range("D3").Value = "2.578.38"
$ws.Range("E3").Value = "  -2.78%  "
$ws.Range("E4").Value = "  +0.09%  "
$ws.Range("D5").Value = "518.82"
$ws.Range("E5").Value = "  -0.68%  "
$ws.Range("D6").Value = "143.31"
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").Value = "0.997"
$ws.Range("E7").Value = "  -0.16%  "
$ws.Range("D8").Value = "0.566"
$ws.Range("E8").Value = "  -1.31%  "
$ws.Range("D9").Value = "2.600.61"
$ws.Range("E9").Value = "  -2.07%  "
$ws.Range("D10").Value = "6.60"
$ws.Range("E10").Value = "  -1.83%  "
$ws.Range("D11").Value = "0.101"
$ws.Range("E11").Value = "  -1.61%  "
$ws.Range("D12").Value = "0.325"
$ws.Range("E12").Value = "  -4.60%  "
$ws.Range("E13").Value = "  -0.78%  "
$ws.Range("D14").Value = "3.033.07"
$ws.Range("E14").Value = "  -2.72%  "
$ws.Range("D15").Value = "57.936.68"
$ws.Range("E15").Value = "  -0.78%  "
$ws.Range("D16").Value = "20.36"
$ws.Range("E16").Value = "  -2.62%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.593.18"
$ws.Range("E17").Value = "  -2.24%  "
$ws.Range("B18").Value = "ShibaInu"
$ws.Range("C18").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D18").Value = "0.0000134"
$ws.Range("E18").Value = "  -1.99%  "
$ws.Range("D19").Value = "339.62"
$ws.Range("E19").Value = "  +0.19%  "
$ws.Range("D20").Value = "4.30"
$ws.Range("E20").Value = "  -2.16%  "
$ws.Range("D21").Value = "10.25"
$ws.Range("E21").Value = "  -1.85%  "
$ws.Range("D22").Value = "6.34"
$ws.Range("E22").Value = "  -0.02%  "
$ws.Range("D23").Value = "0.998"
$ws.Range("E23").Value = "  -0.18%  "
$ws.Range("D24").Value = "65.53"
$ws.Range("E24").Value = "  +1.55%  "
$ws.Range("D25").Value = "0.166"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("B26").Value = "Polygon"
$ws.Range("C26").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D26").Value = "0.403"
$ws.Range("E26").Value = "  -5.35%  "
$ws.Range("B27").Value = "Binance-PegBSC-USD"
$ws.Range("C27").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D27").Value = "0.993"
$ws.Range("E27").Value = "  -0.59%  "
$ws.Range("D28").Value = "2.675.97"
$ws.Range("E28").Value = "  -3.32%  "
$ws.Range("D29").Value = "6.99"
$ws.Range("E29").Value = "  -2.36%  "
$ws.Range("D30").Value = "0.0₃0753"
$ws.Range("E30").Value = "  -5.87%  "
$ws.Range("D31").Value = "0.998"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("D32").Value = "6.36"
$ws.Range("E32").Value = "  -4.74%  "
$ws.Range("E33").Value = "  -0.85%  "
$ws.Range("B34").Value = "EthereumClassic"
$ws.Range("C34").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D34").Value = "18.69"
$ws.Range("E34").Value = "  -1.09%  "
$ws.Range("B35").Value = "Monero"
$ws.Range("C35").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D35").Value = "149.26"
$ws.Range("E35").Value = "  -1.58%  "
$ws.Range("D36").Value = "4.02"
$ws.Range("E36").Value = "  -3.55%  "
$ws.Range("D37").Value = "1.15"
$ws.Range("E37").Value = "  -2.91%  "
$ws.Range("D38").Value = "0.875"
$ws.Range("E38").Value = "  -4.76%  "
$ws.Range("D39").Value = "36.20"
$ws.Range("E39").Value = "  -1.74%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Value = "1.46"
$ws.Range("E40").Value = "  +1.52%  "
$ws.Range("B41").Value = "Fetch.AI"
$ws.Range("C41").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D41").Value = "0.837"
$ws.Range("E41").Value = "  -4.04%  "
$ws.Range("D42").Value = "3.54"
$ws.Range("D43").Value = "0.997"
$ws.Range("E43").Value = "  -0.17%  "
$ws.Range("D44").Value = "271.72"
$ws.Range("E44").Value = "  -1.44%  "
$ws.Range("B45").Value = "Stellar"
$ws.Range("C45").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D45").Value = "0.0955"
$ws.Range("E45").Value = "  -1.88%  "
$ws.Range("B46").Value = "WhiteBITCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D46").Value = "10.66"
$ws.Range("E46").Value = "  +0.37%  "
$ws.Range("D47").Value = "0.590"
$ws.Range("E47").Value = "  -3.21%  "
$ws.Range("D48").Value = "18.88"
$ws.Range("E48").Value = "  -3.17%  "
$ws.Range("D49").Value = "0.0524"
$ws.Range("E49").Value = "  -2.32%  "
$ws.Range("D50").Value = "1.978.78"
$ws.Range("E50").Value = "  -3.37%  "
$ws.Range("D51").Value = "4.62"
$ws.Range("E51").Value = "  -1.66%  "
